$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Arun@506$"
$ws.Range("B3").Value = "Priya@028$"
$ws.Range("B4").Value = "Ramesh@541$"
$ws.Range("B5").Value = "Divya@564$"
$ws.Range("B6").Value = "Karthik@097$"
$ws.Range("B7").Value = "Sneha@010$"
$ws.Range("B8").Value = "Vijay@519$"
$ws.Range("B9").Value = "Meera@826$"
$ws.Range("B10").Value = "Suresh@866$"
$ws.Range("B11").Value = "Arun@955$"
$ws.Range("B12").Value = "Rajesh@971$"
$ws.Range("B13").Value = "Pooja@818$"
$ws.Range("B14").Value = "Senthil@556$"
$ws.Range("B15").Value = "Lakshmi@896$"
$ws.Range("B16").Value = "Sivakumar@743$"
$ws.Range("B17").Value = "Revathi@073$"
$ws.Range("B18").Value = "Gopinath@027$"
$ws.Range("B19").Value = "Shanthi@935$"
$ws.Range("B20").Value = "Balaji@252$"
$ws.Range("B21").Value = "Kavitha@576$"
$ws.Range("B22").Value = "Deepak@878$"
$ws.Range("B23").Value = "Anjali@734$"
$ws.Range("B24").Value = "Manoj@869$"
$ws.Range("B25").Value = "Ritu@692$"
$ws.Range("B26").Value = "Harish@015$"
$ws.Range("B27").Value = "Neha@501$"
$ws.Range("B28").Value = "Prakash@574$"
$ws.Range("B29").Value = "Asha@702$"
$ws.Range("B30").Value = "Sanjay@952$"
$ws.Range("B31").Value = "Vidya@556$"
